$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 66: High-Load Ball Joint Rod End, McMaster, 4482T6, $19.97
$ws.Range("A66").Formula = "=A65+1"
$ws.Range("B66").Value = "High-Load Ball Joint Rod End"
$ws.Range("C66").Value = "McMaster"
$ws.Range("D66").Value = "4482T6"
$ws.Range("E66").Value = 19.97
$ws.Range("E66").NumberFormat = '"$"#,##0.00_);[Red]("$"#,##0.00)'

# Row 67: Front Upper A-Arm, Assembly
$ws.Range("A67").Formula = "=A66+1"
$ws.Range("B67").Value = "Front Upper A-Arm, Assembly"

# Update selection to match post-edit state
$ws.Range("B67").Select()
